# Update trading signals with new market data and timestamps
# Applies new values to "Active Signals", "Summary Dashboard" and
# "Signal History" sheets, and keeps the BUY/SELL cell coloring
# (green for BUY, red for SELL) consistent with the new signal values.

$wb = $excel.ActiveWorkbook

$COLOR_BUY  = 13561798   # BGR for C6EFCE (green)
$COLOR_SELL = 13551615   # BGR for FFC7CE (red)

function Set-SignalFill {
    param($cell, [string]$signal)
    if ($signal -eq "BUY") {
        $cell.Interior.Color = $COLOR_BUY
    } else {
        $cell.Interior.Color = $COLOR_SELL
    }
}

function Set-TextValue {
    # Forces the cell to keep/store a literal text value (not get
    # auto-converted to a number/percentage/date by Excel).
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------
# Sheet 1: "Active Signals"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Active Signals")

$activeRows = @(
    @{ Row=2;  A="2025-07-28 20:09"; B="XAUGBP"; C="SELL"; D=2093.51867; E=2093.52187; F=2093.5146;  G=0.04; H="84.0%"; I=1.27; J="Active" },
    @{ Row=3;  A="2025-07-28 20:33"; B="XAUCAD"; C="SELL"; D=3637.02642; E=3637.02962; F=3637.02181; G=0.05; H="90.0%"; I=1.44; J="Active" },
    @{ Row=4;  A="2025-07-28 20:05"; B="USDCAD"; C="BUY";  D=1.36249;    E=1.35919;    F=1.37219;    G=0.03; H="84.0%"; I=2.94; J="Active" },
    @{ Row=5;  A="2025-07-28 19:54"; B="XAUCHF"; C="BUY";  D=2338.51219; E=2338.50814; F=2338.52196; G=0.01; H="88.0%"; I=2.41; J="Active" },
    @{ Row=6;  A="2025-07-28 20:39"; B="XAUEUR"; C="BUY";  D=2413.91887; E=2413.91527; F=2413.92515; G=0.08; H="83.0%"; I=1.74; J="Active" }
)

foreach ($r in $activeRows) {
    $row = $r.Row
    $ws1.Cells.Item($row, 1).Value = $r.A
    $ws1.Cells.Item($row, 2).Value = $r.B
    $cCell = $ws1.Cells.Item($row, 3)
    $cCell.Value = $r.C
    Set-SignalFill $cCell $r.C
    $ws1.Cells.Item($row, 4).Value = $r.D
    $ws1.Cells.Item($row, 5).Value = $r.E
    $ws1.Cells.Item($row, 6).Value = $r.F
    $ws1.Cells.Item($row, 7).Value = $r.G
    Set-TextValue ($ws1.Cells.Item($row, 8)) $r.H
    $ws1.Cells.Item($row, 9).Value = $r.I
    $ws1.Cells.Item($row, 10).Value = $r.J
}

# ---------------------------------------------------------------
# Sheet 2: "Summary Dashboard"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary Dashboard")

Set-TextValue ($ws2.Range("B7")) "83.8%"
Set-TextValue ($ws2.Range("B8")) "1.98"
Set-TextValue ($ws2.Range("B9")) "2025-07-28 20:23:57"

# ---------------------------------------------------------------
# Sheet 3: "Signal History"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Signal History")

$historyRows = @(
    @{ Row=2;  A="2025-07-28 20:09"; B="XAUGBP"; C="SELL"; D=2093.51867; E=2093.52187; F=2093.5146;  G=0.04; H=0.84; I=1.27; J="Active" },
    @{ Row=3;  A="2025-07-28 20:44"; B="XAUEUR"; C="SELL"; D=2418.08605; E=2418.09006; F=2418.0814;  G=0.09; H=0.89; I=1.16; J="Filled" },
    @{ Row=4;  A="2025-07-28 20:33"; B="XAUCAD"; C="SELL"; D=3637.02642; E=3637.02962; F=3637.02181; G=0.05; H=0.9;  I=1.44; J="Active" },
    @{ Row=5;  A="2025-07-28 20:09"; B="USDCHF"; C="BUY";  D=0.8848;     E=0.88062;    F=0.89168;    G=0.06; H=0.8;  I=1.65; J="Pending" },
    @{ Row=6;  A="2025-07-28 20:05"; B="USDCAD"; C="BUY";  D=1.36249;    E=1.35919;    F=1.37219;    G=0.03; H=0.84; I=2.94; J="Active" },
    @{ Row=7;  A="2025-07-28 20:44"; B="USDJPY"; C="BUY";  D=149.94536;  E=149.5977;   F=150.9003;   G=0.04; H=0.93; I=2.75; J="Filled" },
    @{ Row=8;  A="2025-07-28 20:42"; B="USDCHF"; C="SELL"; D=0.87869;    E=0.88166;    F=0.8714499999999999; G=0.04; H=0.77; I=2.44; J="Filled" },
    @{ Row=9;  A="2025-07-28 20:50"; B="XAUGBP"; C="BUY";  D=2107.12619; E=2107.1232;  F=2107.13299; G=0.03; H=0.87; I=2.28; J="Pending" },
    @{ Row=10; A="2025-07-28 20:24"; B="XAUEUR"; C="SELL"; D=2416.53418; E=2416.5373;  F=2416.53013; G=0.05; H=0.8;  I=1.3;  J="Filled" },
    @{ Row=11; A="2025-07-28 19:54"; B="XAUCHF"; C="BUY";  D=2338.51219; E=2338.50814; F=2338.52196; G=0.01; H=0.88; I=2.41; J="Active" },
    @{ Row=12; A="2025-07-28 19:54"; B="AUDUSD"; C="SELL"; D=0.65846;    E=0.66216;    F=0.6499;     G=0.03; H=0.82; I=2.32; J="Filled" },
    @{ Row=13; A="2025-07-28 20:33"; B="XAUEUR"; C="BUY";  D=2417.12925; E=2417.12517; F=2417.13821; G=0.08; H=0.78; I=2.2;  J="Pending" },
    @{ Row=14; A="2025-07-28 20:29"; B="XAUUSD"; C="SELL"; D=2661.95755; E=2661.95969; F=2661.95199; G=0.05; H=0.84; I=2.6;  J="Filled" },
    @{ Row=15; A="2025-07-28 20:12"; B="XAUCHF"; C="BUY";  D=2348.2017;  E=2348.19708; F=2348.20753; G=0.03; H=0.78; I=$null; J="Pending" },
    @{ Row=16; A="2025-07-28 20:39"; B="XAUEUR"; C="BUY";  D=2413.91887; E=2413.91527; F=2413.92515; G=0.08; H=0.83; I=1.74; J="Active" }
)

foreach ($r in $historyRows) {
    $row = $r.Row
    $ws3.Cells.Item($row, 1).Value = $r.A
    $ws3.Cells.Item($row, 2).Value = $r.B
    $ws3.Cells.Item($row, 3).Value = $r.C
    $ws3.Cells.Item($row, 4).Value = $r.D
    $ws3.Cells.Item($row, 5).Value = $r.E
    $ws3.Cells.Item($row, 6).Value = $r.F
    $ws3.Cells.Item($row, 7).Value = $r.G
    $ws3.Cells.Item($row, 8).Value = $r.H
    if ($null -ne $r.I) {
        $ws3.Cells.Item($row, 9).Value = $r.I
    }
    $ws3.Cells.Item($row, 10).Value = $r.J
}
